$wb = $excel.ActiveWorkbook

# ---- Sheet: lusid_holdings (new UkEquityActive rows 11-15) ----
$ws1 = $wb.Worksheets.Item("lusid_holdings")

$ws1.Rows.Item(2).Copy()
$ws1.Rows.Item(11).Insert()
$ws1.Rows.Item(2).Copy()
$ws1.Rows.Item(12).Insert()
$ws1.Rows.Item(2).Copy()
$ws1.Rows.Item(13).Insert()
$ws1.Rows.Item(2).Copy()
$ws1.Rows.Item(14).Insert()
$ws1.Rows.Item(2).Copy()
$ws1.Rows.Item(15).Insert()

$ws1.Range("A11:A15").Value = "UkEquityActive"

$ws1.Range("C11").Value = "GB0031348658"
$ws1.Range("D11").Value = 1000

$ws1.Range("C12").Value = "GB00BH0P3Z91"
$ws1.Range("D12").Value = 2000

$ws1.Range("C13").Value = "GB0031743007"
$ws1.Range("D13").Value = 1500

$ws1.Range("C14").Value = "GB0009252882"
$ws1.Range("D14").Value = 1000

$ws1.Range("C15").Value = "JE00B4T3BW64"
$ws1.Range("D15").Value = 1750

$ws1.Columns.Item(4).ColumnWidth = 3.8

# ---- Sheet: ext_holdings (mirror of the same new rows) ----
$ws3 = $wb.Worksheets.Item("ext_holdings")

$ws3.Rows.Item(2).Copy()
$ws3.Rows.Item(11).Insert()
$ws3.Rows.Item(2).Copy()
$ws3.Rows.Item(12).Insert()
$ws3.Rows.Item(2).Copy()
$ws3.Rows.Item(13).Insert()
$ws3.Rows.Item(2).Copy()
$ws3.Rows.Item(14).Insert()
$ws3.Rows.Item(2).Copy()
$ws3.Rows.Item(15).Insert()

$ws3.Range("A11:A15").Value = "UkEquityActive"

$ws3.Range("C11").Value = "GB0031348658"
$ws3.Range("D11").Value = 1000

$ws3.Range("C12").Value = "GB00BH0P3Z91"
$ws3.Range("D12").Value = 2000

$ws3.Range("C13").Value = "GB0031743007"
$ws3.Range("D13").Value = 1500

$ws3.Range("C14").Value = "GB0009252882"
$ws3.Range("D14").Value = 1000

$ws3.Range("C15").Value = "JE00B4T3BW64"
$ws3.Range("D15").Value = 1750

$ws3.Columns.Item(1).ColumnWidth = 12.9

# ---- Selections / active sheet to match the final view ----
$ws1.Range("A11:E15").Select()
$ws3.Range("J9").Select()
